$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 37, shifting existing rows 37-67 down to 38-68.
$ws.Rows("37:37").Insert()

# Populate the newly inserted row 37 with its data.
$ws.Cells.Item(37, 1).Value = 7
$ws.Cells.Item(37, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(37, 3).Value = "Ñuble"
$ws.Cells.Item(37, 4).Value = 44596
$ws.Cells.Item(37, 5).Value = 16
$ws.Cells.Item(37, 6).Value = 100112031
$ws.Cells.Item(37, 7).Value = "Poroto verde"
$ws.Cells.Item(37, 8).Value = "Sin especificar"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 29000
$ws.Cells.Item(37, 12).Value = 30000
$ws.Cells.Item(37, 13).Value = 29500
$ws.Cells.Item(37, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(37, 15).Value = "Región del Maule"
$ws.Cells.Item(37, 16).Value = 1180
$ws.Cells.Item(37, 17).Value = 25
$ws.Cells.Item(37, 18).Value = "Hortaliza"
